$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.99000000000062
$ws.Range("H2").Value = 0.0001917361827612485
$ws.Range("I2").Value = 0.0001917361827612485
$ws.Range("L2").Value = 50.93300869671616
$ws.Range("M2").Value = "[19.549681799370816, 82.3163355940615]"
$ws.Range("N2").Value = 0.002073918683734632
$ws.Range("O2").Value = 0.002073918683734632
$ws.Range("P2").Value = 2.056658253701427
$ws.Range("Q2").Value = "[1.4528686746331179, 2.6604478327697354]"
$ws.Range("R2").Value = [double]"1.641021252218877e-08"
$ws.Range("S2").Value = [double]"1.641021252218877e-08"
$ws.Range("T2").Value = 63.29660530194285
$ws.Range("U2").Value = "[46.720280901567776, 79.87292970231792]"
$ws.Range("V2").Value = [double]"9.7506447360729e-10"
$ws.Range("W2").Value = [double]"9.7506447360729e-10"
$ws.Range("X2").Value = 17.48276276276318
$ws.Range("Y2").Value = 14.98522522522558
$ws.Range("Z2").Value = 19.98030030030078

# Row 3
$ws.Range("F3").Value = 25.99000000000062
$ws.Range("H3").Value = [double]"3.5009012053e-05"
$ws.Range("I3").Value = [double]"3.5009012053e-05"
$ws.Range("L3").Value = 58.11955522168175
$ws.Range("M3").Value = "[28.92021397864555, 87.31889646471795]"
$ws.Range("N3").Value = 0.0002268584270908036
$ws.Range("O3").Value = 0.0002268584270908036
$ws.Range("P3").Value = 1.327079178993888
$ws.Range("Q3").Value = "[0.7232895999255788, 1.9308687580621964]"
$ws.Range("R3").Value = [double]"6.029597441026269e-05"
$ws.Range("S3").Value = [double]"6.029597441026269e-05"
$ws.Range("T3").Value = 61.41434507358598
$ws.Range("U3").Value = "[44.80468157602223, 78.02400857114972]"
$ws.Range("V3").Value = [double]"2.223980333226905e-09"
$ws.Range("W3").Value = [double]"2.223980333226905e-09"
$ws.Range("X3").Value = 20.50062062062111
$ws.Range("Y3").Value = 18.00308308308351
$ws.Range("Z3").Value = 22.99815815815871

# Row 4
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 25.99000000000062
$ws.Range("H4").Value = [double]"3.127127655710993e-07"
$ws.Range("I4").Value = [double]"3.127127655710993e-07"
$ws.Range("L4").Value = 53.03690071003057
$ws.Range("M4").Value = "[30.846534427878552, 75.22726699218259]"
$ws.Range("N4").Value = [double]"1.705007434504857e-05"
$ws.Range("O4").Value = [double]"1.705007434504857e-05"
$ws.Range("P4").Value = 0.7610264486173479
$ws.Range("Q4").Value = "[0.37107901213573236, 1.1509738850989635]"
$ws.Range("R4").Value = 0.0002891508033022738
$ws.Range("S4").Value = 0.0002891508033022738
$ws.Range("T4").Value = 53.2694087833953
$ws.Range("U4").Value = "[41.9641588292026, 64.574658737588]"
$ws.Range("V4").Value = [double]"2.608802063264193e-12"
$ws.Range("W4").Value = [double]"2.608802063264193e-12"
$ws.Range("X4").Value = 22.84206206206261
$ws.Range("Y4").Value = 21.22906906906958
$ws.Range("Z4").Value = 24.45505505505564

# Row 5
$ws.Range("F5").Value = 25.99000000000062
$ws.Range("H5").Value = [double]"5.545352655333957e-05"
$ws.Range("I5").Value = [double]"5.545352655333957e-05"
$ws.Range("L5").Value = 56.44680405018264
$ws.Range("M5").Value = "[28.895101740508323, 83.99850635985695]"
$ws.Range("N5").Value = 0.0001570482102462023
$ws.Range("O5").Value = 0.0001570482102462023
$ws.Range("P5").Value = 0.1823947686768852
$ws.Range("Q5").Value = "[-0.3710790121357306, 0.735868549489501]"
$ws.Range("R5").Value = 0.5102440563537138
$ws.Range("S5").Value = 0.5102440563537138
$ws.Range("T5").Value = 78.2070842637626
$ws.Range("U5").Value = "[62.65536400248496, 93.75880452504023]"
$ws.Range("V5").Value = [double]"3.479438959175241e-13"
$ws.Range("W5").Value = [double]"3.479438959175241e-13"
$ws.Range("X5").Value = 25.23553553553614
$ws.Range("Y5").Value = 22.94612612612668
$ws.Range("Z5").Value = 27.5249449449456

# Row 6
$ws.Range("F6").Value = 25.99000000000062
$ws.Range("H6").Value = 0.000984022428868192
$ws.Range("I6").Value = 0.000984022428868192
$ws.Range("L6").Value = 46.22644032048697
$ws.Range("M6").Value = "[19.944762351083327, 72.50811828989062]"
$ws.Range("N6").Value = 0.0009361800686653599
$ws.Range("O6").Value = 0.0009361800686653599
$ws.Range("P6").Value = -0.08805264694746207
$ws.Range("Q6").Value = "[-0.7987632973091179, 0.6226580034141938]"
$ws.Range("R6").Value = 0.8040825945181713
$ws.Range("S6").Value = 0.8040825945181713
$ws.Range("T6").Value = 61.8356032817584
$ws.Range("U6").Value = "[45.95237255791106, 77.71883400560573]"
$ws.Range("V6").Value = [double]"5.875233632934851e-10"
$ws.Range("W6").Value = [double]"5.875233632934851e-10"
$ws.Range("X6").Value = 0.364224224224234
$ws.Range("Y6").Value = -2.57558558558565
$ws.Range("Z6").Value = 3.304034034034118

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 25.99000000000062
$ws.Range("H7").Value = [double]"8.482601002723911e-05"
$ws.Range("I7").Value = [double]"8.482601002723911e-05"
$ws.Range("L7").Value = 52.13179424583979
$ws.Range("M7").Value = "[22.16279319266613, 82.10079529901346]"
$ws.Range("N7").Value = 0.00105034959879835
$ws.Range("O7").Value = 0.00105034959879835
$ws.Range("P7").Value = -0.5031579825569237
$ws.Range("Q7").Value = "[-1.0692107129334634, 0.06289474781961601]"
$ws.Range("R7").Value = 0.08013264357139449
$ws.Range("S7").Value = 0.08013264357139449
$ws.Range("T7").Value = 56.46871175451833
$ws.Range("U7").Value = "[41.002328490674415, 71.93509501836225]"
$ws.Range("V7").Value = [double]"3.055205644031389e-09"
$ws.Range("W7").Value = [double]"3.055205644031389e-09"
$ws.Range("X7").Value = 2.081281281281331
$ws.Range("Y7").Value = -0.2601601601601682
$ws.Range("Z7").Value = 4.422722722722831

# Row 8
$ws.Range("F8").Value = 25.99000000000062
$ws.Range("H8").Value = [double]"2.063498254489815e-07"
$ws.Range("I8").Value = [double]"2.063498254489815e-07"
$ws.Range("L8").Value = 47.72349955150602
$ws.Range("M8").Value = "[28.35218886301439, 67.09481023999766]"
$ws.Range("N8").Value = [double]"1.043948896839453e-05"
$ws.Range("O8").Value = [double]"1.043948896839453e-05"
$ws.Range("P8").Value = -1.25789495639231
$ws.Range("Q8").Value = "[-1.6981581911296182, -0.8176317216550011]"
$ws.Range("R8").Value = [double]"7.214686643663981e-07"
$ws.Range("S8").Value = [double]"7.214686643663981e-07"
$ws.Range("T8").Value = 52.57426220678582
$ws.Range("U8").Value = "[41.87330236465812, 63.27522204891351]"
$ws.Range("V8").Value = [double]"7.223110998211268e-13"
$ws.Range("W8").Value = [double]"7.223110998211268e-13"
$ws.Range("X8").Value = 5.20320320320333
$ws.Range("Y8").Value = 3.382082082082164
$ws.Range("Z8").Value = 7.024324324324496

# Row 9
$ws.Range("F9").Value = 25.68000000000058
$ws.Range("H9").Value = [double]"2.583200835570842e-06"
$ws.Range("I9").Value = [double]"2.583200835570842e-06"
$ws.Range("L9").Value = 60.07767044658139
$ws.Range("M9").Value = "[31.354985862680596, 88.80035503048218]"
$ws.Range("N9").Value = 0.0001195414907249059
$ws.Range("O9").Value = 0.0001195414907249059
$ws.Range("P9").Value = -0.704421175579693
$ws.Range("Q9").Value = "[-1.1698423094448476, -0.23900004171453837]"
$ws.Range("R9").Value = 0.003844672146208872
$ws.Range("S9").Value = 0.003844672146208872
$ws.Range("T9").Value = 65.85651588661514
$ws.Range("U9").Value = "[51.10271214723895, 80.61031962599134]"
$ws.Range("V9").Value = [double]"1.307509656101047e-11"
$ws.Range("W9").Value = [double]"1.307509656101047e-11"
$ws.Range("X9").Value = 2.8790390390391
$ws.Range("Y9").Value = 0.9768168168168345
$ws.Range("Z9").Value = 4.781261261261365

# Row 10
$ws.Range("F10").Value = 25.68000000000058
$ws.Range("H10").Value = 0.0009530306428154889
$ws.Range("I10").Value = 0.0009530306428154889
$ws.Range("L10").Value = 42.12068894969972
$ws.Range("M10").Value = "[18.101524740630467, 66.13985315876897]"
$ws.Range("N10").Value = 0.000965943272999592
$ws.Range("O10").Value = 0.000965943272999592
$ws.Range("P10").Value = -1.622684493746079
$ws.Range("Q10").Value = "[-2.427737265837157, -0.8176317216550011]"
$ws.Range("R10").Value = 0.00019363775407788
$ws.Range("S10").Value = 0.00019363775407788
$ws.Range("T10").Value = 63.50333393952781
$ws.Range("U10").Value = "[48.18388043742501, 78.82278744163062]"
$ws.Range("V10").Value = [double]"1.076201350258543e-10"
$ws.Range("W10").Value = [double]"1.076201350258543e-10"
$ws.Range("X10").Value = 6.632072072072219
$ws.Range("Y10").Value = 3.341741741741815
$ws.Range("Z10").Value = 9.922402402402623

# Row 11
$ws.Range("F11").Value = 25.68000000000058
$ws.Range("H11").Value = 0.001292858091501148
$ws.Range("I11").Value = 0.001292858091501148
$ws.Range("L11").Value = 42.88256448256066
$ws.Range("M11").Value = "[14.064511999093341, 71.70061696602798]"
$ws.Range("N11").Value = 0.00442469251483768
$ws.Range("O11").Value = 0.00442469251483768
$ws.Range("P11").Value = -1.245316006828387
$ws.Range("Q11").Value = "[-2.0000529806637726, -0.49057903299300065]"
$ws.Range("R11").Value = 0.001774458769027376
$ws.Range("S11").Value = 0.001774458769027376
$ws.Range("T11").Value = 66.18269076825302
$ws.Range("U11").Value = "[50.05917517966361, 82.30620635684244]"
$ws.Range("V11").Value = [double]"1.411708527854216e-10"
$ws.Range("W11").Value = [double]"1.411708527854216e-10"
$ws.Range("X11").Value = 5.089729729729843
$ws.Range("Y11").Value = 2.005045045045089
$ws.Range("Z11").Value = 8.174414414414597
